$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strain TDY1452 (rows 8-10) marker_1 changed from NAT to G418
$ws.Range("J8").Value = "G418"
$ws.Range("J9").Value = "G418"
$ws.Range("J10").Value = "G418"

# Reflect the active cell selection recorded in the saved workbook
$ws.Range("P7").Select()
